# The "all_finals_untoned" lookup table had a spurious row pairing
# med="$eh" with rime/py="ê" (a mistaken duplicate of the "e"/"eh" rows).
# Remove it; Excel reflows the remaining rows up and drops the now-unused
# "ê" shared string automatically on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_finals_untoned")
$ws.Activate()

$ws.Rows.Item(5).Delete()

# Leave the view where the user ended up after the edit.
$ws.Range("I32").Select()
